$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.767.85"
$ws.Range("E2").Value = "'  -1.84%  "
$ws.Range("D3").Value = "'2.369.34"
$ws.Range("E3").Value = "'  -2.14%  "
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("D5").Value = "'558.89"
$ws.Range("E5").Value = "'  -1.82%  "
$ws.Range("D6").Value = "'137.60"
$ws.Range("E6").Value = "'  -1.29%  "
$ws.Range("E7").Value = "'  -0.26%  "
$ws.Range("D8").Value = "'0.529"
$ws.Range("E8").Value = "'  +0.73%  "
$ws.Range("D9").Value = "'2.364.91"
$ws.Range("E9").Value = "'  -1.58%  "
$ws.Range("E10").Value = "'  -2.23%  "
$ws.Range("E11").Value = "'  -1.00%  "
$ws.Range("E12").Value = "'  +0.77%  "
$ws.Range("D13").Value = "'0.337"
$ws.Range("E13").Value = "'  +0.30%  "
$ws.Range("D14").Value = "'25.56"
$ws.Range("E14").Value = "'  -1.24%  "
$ws.Range("D15").Value = "'2.794.38"
$ws.Range("E15").Value = "'  -0.91%  "
$ws.Range("E16").Value = "'  -2.47%  "
$ws.Range("D17").Value = "'59.645.15"
$ws.Range("E17").Value = "'  -1.97%  "
$ws.Range("D18").Value = "'2.356.79"
$ws.Range("E18").Value = "'  -2.18%  "
$ws.Range("D19").Value = "'8.01"
$ws.Range("E19").Value = "'  +11.03%  "
$ws.Range("E20").Value = "'  -0.26%  "
$ws.Range("D21").Value = "'321.30"
$ws.Range("E21").Value = "'  +0.31%  "
$ws.Range("E22").Value = "'  +1.08%  "
$ws.Range("E23").Value = "'  -1.62%  "
$ws.Range("D25").Value = "'1.81"
$ws.Range("E25").Value = "'  -3.36%  "
$ws.Range("D26").Value = "'64.17"
$ws.Range("E26").Value = "'  -0.91%  "
$ws.Range("D27").Value = "'558.34"
$ws.Range("E27").Value = "'  -2.19%  "
$ws.Range("E28").Value = "'  -6.33%  "
$ws.Range("D30").Value = "'0.0₃0921"
$ws.Range("E30").Value = "'  +1.92%  "
$ws.Range("D31").Value = "'7.99"
$ws.Range("E31").Value = "'  +2.38%  "
$ws.Range("E32").Value = "'  -2.00%  "
$ws.Range("E33").Value = "'  -2.71%  "
$ws.Range("E34").Value = "'  -1.75%  "
$ws.Range("E35").Value = "'  -0.87%  "
$ws.Range("E36").Value = "'  +3.75%  "
$ws.Range("D37").Value = "'152.92"
$ws.Range("E37").Value = "'  +2.54%  "
$ws.Range("D38").Value = "'0.366"
$ws.Range("E38").Value = "'  +0.20%  "
$ws.Range("E39").Value = "'  -1.20%  "
$ws.Range("D40").Value = "'18.15"
$ws.Range("E40").Value = "'  +0.54%  "
$ws.Range("E41").Value = "'  -1.12%  "
$ws.Range("E42").Value = "'  -0.02%  "
$ws.Range("D43").Value = "'41.39"
$ws.Range("E43").Value = "'  -0.79%  "
$ws.Range("E44").Value = "'  -0.74%  "
$ws.Range("E45").Value = "'  +3.94%  "
$ws.Range("E46").Value = "'  +6.08%  "
$ws.Range("D47").Value = "'138.38"
$ws.Range("E47").Value = "'  -0.94%  "
$ws.Range("E48").Value = "'  +0.91%  "
$ws.Range("E49").Value = "'  -0.84%  "
$ws.Range("D50").Value = "'0.0499"
$ws.Range("E50").Value = "'  -0.55%  "
$ws.Range("E51").Value = "'  -0.73%  "
